$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2223.7322
$ws.Cells.Item(15, 9).Value = 2223.7322
$ws.Cells.Item(15, 11).Value = 6671.196599999999
$ws.Cells.Item(15, 13).Value = -6502.196599999999

$ws.Cells.Item(43, 8).Value = 8798.700000000001
$ws.Cells.Item(43, 9).Value = 8799.014999999999
$ws.Cells.Item(43, 10).Value = 8798
$ws.Cells.Item(43, 11).Value = 8799.014999999999
$ws.Cells.Item(43, 12).Value = 8798
$ws.Cells.Item(43, 13).Value = -8730.014999999999
$ws.Cells.Item(43, 14).Value = -8936

$ws.Cells.Item(62, 8).Value = 4146.9565
$ws.Cells.Item(62, 9).Value = 3140.8235
$ws.Cells.Item(62, 10).Value = 6997.6665
$ws.Cells.Item(62, 11).Value = 3140.8235
$ws.Cells.Item(62, 12).Value = 6997.6665
$ws.Cells.Item(62, 13).Value = -2516.8235
$ws.Cells.Item(62, 14).Value = -8245.666499999999

$ws.Cells.Item(65, 8).Value = 4146.9565
$ws.Cells.Item(65, 9).Value = 3140.8235
$ws.Cells.Item(65, 10).Value = 6997.6665
$ws.Cells.Item(65, 11).Value = 15704.1175
$ws.Cells.Item(65, 12).Value = 34988.3325
$ws.Cells.Item(65, 13).Value = -12584.1175
$ws.Cells.Item(65, 14).Value = -41228.3325

$ws.Cells.Item(69, 8).Value = 6207.0435
$ws.Cells.Item(69, 9).Value = 4290
$ws.Cells.Item(69, 10).Value = 7681.6924
$ws.Cells.Item(69, 11).Value = 12870
$ws.Cells.Item(69, 12).Value = 23045.0772
$ws.Cells.Item(69, 13).Value = -11996
$ws.Cells.Item(69, 14).Value = -24793.0772

$ws.Cells.Item(72, 8).Value = 6207.0435
$ws.Cells.Item(72, 9).Value = 4290
$ws.Cells.Item(72, 10).Value = 7681.6924
$ws.Cells.Item(72, 11).Value = 38610
$ws.Cells.Item(72, 12).Value = 69135.2316
$ws.Cells.Item(72, 13).Value = -34242
$ws.Cells.Item(72, 14).Value = -77871.2316

$ws.Cells.Item(116, 8).Value = 3368
$ws.Cells.Item(116, 9).Value = 2924
$ws.Cells.Item(116, 11).Value = 2924
$ws.Cells.Item(116, 13).Value = 518

$ws.Cells.Item(129, 8).Value = 1964.9445
$ws.Cells.Item(129, 9).Value = 935
$ws.Cells.Item(129, 11).Value = 2805
$ws.Cells.Item(129, 13).Value = 2195

$ws.Cells.Item(137, 8).Value = 6945.6733
$ws.Cells.Item(137, 9).Value = 3185.2693
$ws.Cells.Item(137, 10).Value = 11196.565
$ws.Cells.Item(137, 11).Value = 9555.8079
$ws.Cells.Item(137, 12).Value = 33589.695
$ws.Cells.Item(137, 13).Value = -7005.8079
$ws.Cells.Item(137, 14).Value = -38689.695

$ws.Cells.Item(138, 8).Value = 2611.8455
$ws.Cells.Item(138, 9).Value = 2243.5
$ws.Cells.Item(138, 10).Value = 2663.8472
$ws.Cells.Item(138, 11).Value = 6730.5
$ws.Cells.Item(138, 12).Value = 7991.5416
$ws.Cells.Item(138, 13).Value = -1590.5
$ws.Cells.Item(138, 14).Value = -18271.5416

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6373.14
$ws.Cells.Item(32, 9).Value = 1116.0133
$ws.Cells.Item(32, 10).Value = 22144.52
$ws.Cells.Item(32, 11).Value = 1116.0133
$ws.Cells.Item(32, 12).Value = 22144.52
$ws.Cells.Item(32, 13).Value = -829.0133000000001
$ws.Cells.Item(32, 14).Value = -22718.52

$ws.Cells.Item(45, 8).Value = 8987.5
$ws.Cells.Item(45, 9).Value = 8987.5
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 8987.5
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -8610.5
$ws.Cells.Item(45, 14).ClearContents()

$ws.Cells.Item(74, 8).Value = 18659.48
$ws.Cells.Item(74, 9).Value = 2584.1875
$ws.Cells.Item(74, 10).Value = 47237.777
$ws.Cells.Item(74, 11).Value = 2584.1875
$ws.Cells.Item(74, 12).Value = 47237.777
$ws.Cells.Item(74, 13).Value = -1710.1875
$ws.Cells.Item(74, 14).Value = -48985.777

$ws.Cells.Item(77, 8).Value = 18659.48
$ws.Cells.Item(77, 9).Value = 2584.1875
$ws.Cells.Item(77, 10).Value = 47237.777
$ws.Cells.Item(77, 11).Value = 12920.9375
$ws.Cells.Item(77, 12).Value = 236188.885
$ws.Cells.Item(77, 13).Value = -8552.9375
$ws.Cells.Item(77, 14).Value = -244924.885

$ws.Cells.Item(122, 8).Value = 3010.8572
$ws.Cells.Item(122, 9).Value = 1536.125
$ws.Cells.Item(122, 10).Value = 6228.4546
$ws.Cells.Item(122, 11).Value = 4608.375
$ws.Cells.Item(122, 12).Value = 18685.3638
$ws.Cells.Item(122, 13).Value = -2158.375
$ws.Cells.Item(122, 14).Value = -23585.3638

$ws.Cells.Item(132, 8).Value = 1793773.8
$ws.Cells.Item(132, 9).Value = 2657.6743
$ws.Cells.Item(132, 10).Value = 7718234.5
$ws.Cells.Item(132, 11).Value = 7973.0229
$ws.Cells.Item(132, 12).Value = 23154703.5
$ws.Cells.Item(132, 13).Value = -5443.0229
$ws.Cells.Item(132, 14).Value = -23159763.5

$ws.Cells.Item(133, 8).Value = 61249.5
$ws.Cells.Item(133, 10).Value = 61249.5
$ws.Cells.Item(133, 12).Value = 61249.5
$ws.Cells.Item(133, 14).Value = -66309.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).ClearContents()

$ws.Cells.Item(31, 8).Value = 15581.807
$ws.Cells.Item(31, 9).Value = 7148.263
$ws.Cells.Item(31, 10).Value = 28934.916
$ws.Cells.Item(31, 11).Value = 7148.263
$ws.Cells.Item(31, 12).Value = 28934.916
$ws.Cells.Item(31, 13).Value = -6853.263
$ws.Cells.Item(31, 14).Value = -29524.916

$ws.Cells.Item(34, 8).Value = 15581.807
$ws.Cells.Item(34, 9).Value = 7148.263
$ws.Cells.Item(34, 10).Value = 28934.916
$ws.Cells.Item(34, 11).Value = 7148.263
$ws.Cells.Item(34, 12).Value = 28934.916
$ws.Cells.Item(34, 13).Value = -6946.263
$ws.Cells.Item(34, 14).Value = -29338.916

$ws.Cells.Item(125, 8).Value = 76675
$ws.Cells.Item(125, 10).Value = 76675
$ws.Cells.Item(125, 12).Value = 76675
$ws.Cells.Item(125, 14).Value = -81595

$ws.Cells.Item(141, 8).Value = 108719.445
$ws.Cells.Item(141, 10).Value = 108719.445
$ws.Cells.Item(141, 12).Value = 108719.445
$ws.Cells.Item(141, 14).Value = -119079.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 20000690
$ws.Cells.Item(32, 9).Value = 50000224
$ws.Cells.Item(32, 10).Value = 12500807
$ws.Cells.Item(32, 11).Value = 150000672
$ws.Cells.Item(32, 12).Value = 37502421
$ws.Cells.Item(32, 13).Value = -150000389
$ws.Cells.Item(32, 14).Value = -37502987

$ws.Cells.Item(46, 8).Value = 5498.923
$ws.Cells.Item(46, 9).Value = 520.1
$ws.Cells.Item(46, 10).Value = 22095
$ws.Cells.Item(46, 11).Value = 1560.3
$ws.Cells.Item(46, 12).Value = 66285
$ws.Cells.Item(46, 13).Value = -1469.3
$ws.Cells.Item(46, 14).Value = -66467

$ws.Cells.Item(139, 8).Value = 12256.9
$ws.Cells.Item(139, 9).Value = 47712.668
$ws.Cells.Item(139, 10).Value = 6000
$ws.Cells.Item(139, 11).Value = 143138.004
$ws.Cells.Item(139, 12).Value = 18000
$ws.Cells.Item(139, 13).Value = -137998.004
$ws.Cells.Item(139, 14).Value = -28280

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3976896.5
$ws.Cells.Item(126, 9).Value = 12784.2
$ws.Cells.Item(126, 10).Value = 5501555
$ws.Cells.Item(126, 11).Value = 38352.60000000001
$ws.Cells.Item(126, 12).Value = 16504665
$ws.Cells.Item(126, 13).Value = -35882.60000000001
$ws.Cells.Item(126, 14).Value = -16509605

$ws.Cells.Item(141, 8).Value = 122497.5
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 122497.5
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 122497.5
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -132857.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2865.9524
$ws.Cells.Item(22, 9).Value = 2915.6667
$ws.Cells.Item(22, 10).Value = 2799.6667
$ws.Cells.Item(22, 11).Value = 2915.6667
$ws.Cells.Item(22, 12).Value = 2799.6667
$ws.Cells.Item(22, 13).Value = -2620.6667
$ws.Cells.Item(22, 14).Value = -3389.6667

$ws.Cells.Item(27, 8).Value = 2865.9524
$ws.Cells.Item(27, 9).Value = 2915.6667
$ws.Cells.Item(27, 10).Value = 2799.6667
$ws.Cells.Item(27, 11).Value = 2915.6667
$ws.Cells.Item(27, 12).Value = 2799.6667
$ws.Cells.Item(27, 13).Value = -2808.6667
$ws.Cells.Item(27, 14).Value = -3013.6667

$ws.Cells.Item(46, 8).Value = 2859.4614
$ws.Cells.Item(46, 9).Value = 1647.125
$ws.Cells.Item(46, 11).Value = 1647.125
$ws.Cells.Item(46, 13).Value = -1459.125

$ws.Cells.Item(55, 8).Value = 1164.8286
$ws.Cells.Item(55, 9).Value = 777.04
$ws.Cells.Item(55, 10).Value = 2134.3
$ws.Cells.Item(55, 11).Value = 777.04
$ws.Cells.Item(55, 12).Value = 2134.3
$ws.Cells.Item(55, 13).Value = -604.04
$ws.Cells.Item(55, 14).Value = -2480.3

$ws.Cells.Item(96, 8).Value = 14498.75
$ws.Cells.Item(96, 10).Value = 14498.75
$ws.Cells.Item(96, 12).Value = 14498.75
$ws.Cells.Item(96, 14).Value = -19990.75

$ws.Cells.Item(132, 8).Value = 1491968.4
$ws.Cells.Item(132, 9).Value = 2369.25
$ws.Cells.Item(132, 10).Value = 3658658
$ws.Cells.Item(132, 11).Value = 7107.75
$ws.Cells.Item(132, 12).Value = 10975974
$ws.Cells.Item(132, 13).Value = -4577.75
$ws.Cells.Item(132, 14).Value = -10981034

$ws.Cells.Item(136, 8).Value = 14713.697
$ws.Cells.Item(136, 9).Value = 12678.053
$ws.Cells.Item(136, 10).Value = 17476.357
$ws.Cells.Item(136, 11).Value = 38034.159
$ws.Cells.Item(136, 12).Value = 52429.071
$ws.Cells.Item(136, 13).Value = -35484.159
$ws.Cells.Item(136, 14).Value = -57529.071

$ws.Cells.Item(139, 8).Value = 59786.25
$ws.Cells.Item(139, 10).Value = 69715
$ws.Cells.Item(139, 12).Value = 69715
$ws.Cells.Item(139, 14).Value = -79995

$ws.Cells.Item(140, 8).Value = 177911.33
$ws.Cells.Item(140, 10).Value = 177911.33
$ws.Cells.Item(140, 12).Value = 177911.33
$ws.Cells.Item(140, 14).Value = -188271.33

$ws.Cells.Item(141, 8).Value = 68999.8
$ws.Cells.Item(141, 10).Value = 68999.8
$ws.Cells.Item(141, 12).Value = 68999.8
$ws.Cells.Item(141, 14).Value = -79359.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2564.4285
$ws.Cells.Item(96, 9).Value = 1998
$ws.Cells.Item(96, 10).Value = 2658.8333
$ws.Cells.Item(96, 11).Value = 1998
$ws.Cells.Item(96, 12).Value = 2658.8333
$ws.Cells.Item(96, 13).Value = -625
$ws.Cells.Item(96, 14).Value = -5404.8333

$ws.Cells.Item(113, 8).Value = 925.7083
$ws.Cells.Item(113, 9).Value = 558.6429000000001
$ws.Cells.Item(113, 10).Value = 1439.6
$ws.Cells.Item(113, 11).Value = 1675.9287
$ws.Cells.Item(113, 12).Value = 4318.799999999999
$ws.Cells.Item(113, 13).Value = 494.0712999999998
$ws.Cells.Item(113, 14).Value = -8658.799999999999

$ws.Cells.Item(122, 8).Value = 2772.3333
$ws.Cells.Item(122, 9).Value = 1097.5294
$ws.Cells.Item(122, 10).Value = 5619.5
$ws.Cells.Item(122, 11).Value = 3292.5882
$ws.Cells.Item(122, 12).Value = 16858.5
$ws.Cells.Item(122, 13).Value = -842.5881999999997
$ws.Cells.Item(122, 14).Value = -21758.5

$ws.Cells.Item(126, 8).Value = 44633.668
$ws.Cells.Item(126, 9).Value = 51499.875
$ws.Cells.Item(126, 10).Value = 30901.25
$ws.Cells.Item(126, 11).Value = 154499.625
$ws.Cells.Item(126, 12).Value = 92703.75
$ws.Cells.Item(126, 13).Value = -152029.625
$ws.Cells.Item(126, 14).Value = -97643.75

$ws.Cells.Item(132, 8).Value = 5188.7163
$ws.Cells.Item(132, 9).Value = 2758.6597
$ws.Cells.Item(132, 10).Value = 10899.35
$ws.Cells.Item(132, 11).Value = 8275.9791
$ws.Cells.Item(132, 12).Value = 32698.05
$ws.Cells.Item(132, 13).Value = -5745.9791
$ws.Cells.Item(132, 14).Value = -37758.05

$ws.Cells.Item(136, 8).Value = 4944.2295
$ws.Cells.Item(136, 9).Value = 794.51166
$ws.Cells.Item(136, 10).Value = 14857.444
$ws.Cells.Item(136, 11).Value = 2383.53498
$ws.Cells.Item(136, 12).Value = 44572.33199999999
$ws.Cells.Item(136, 13).Value = 166.4650200000001
$ws.Cells.Item(136, 14).Value = -49672.33199999999

$ws.Cells.Item(139, 8).Value = 97000
$ws.Cells.Item(139, 10).Value = 97000
$ws.Cells.Item(139, 12).Value = 97000
$ws.Cells.Item(139, 14).Value = -107280

$ws.Cells.Item(140, 8).Value = 142005.2
$ws.Cells.Item(140, 10).Value = 142005.2
$ws.Cells.Item(140, 12).Value = 142005.2
$ws.Cells.Item(140, 14).Value = -152365.2

$ws.Cells.Item(141, 8).Value = 75200
$ws.Cells.Item(141, 10).Value = 75200
$ws.Cells.Item(141, 12).Value = 75200
$ws.Cells.Item(141, 14).Value = -85560
